$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# "User story" heading gains a trailing space
Replace-Text "User story" "User story "

# User-story paragraphs: append reference numbers / fix typos
Replace-Text `
    "Jako Patrycja chce mieć możliwość wysłania formularza, aby zostać wolontariuszką podczas wydarzenia." `
    "Jako Patrycja chce mieć możliwość wysłania formularza, aby zostać wolontariuszką podczas wydarzenia. [20]"

Replace-Text `
    "Jako Adam chce mieć możliwość rejestracji w systemie, być na bieżąco z informacjami o wydarzeniu oraz móc śledzić live stream z wydarzenia." `
    "Jako Adam chce mieć możliwość rejestracji w systemie, być na bieżąco z informacjami o wydarzeniu oraz móc śledzić live stream z wydarzenia. [2] "

Replace-Text `
    "Jako Użytkownik chce mieć możliwość zalogowania się do systemu, żeby moc korzystac z funkcjonalnosci aplikacji" `
    "Jako Użytkownik chce mieć możliwość zalogowania się do systemu, żeby moc korzystac z funkcjonalnosci aplikacji [5]"

Replace-Text `
    "Jako Damian chce mieć możliwość zakupu biletu na stronie internetowej, by móc uczestniczyć w wydarzeniu." `
    "Jako Damian chce mieć możliwość zakupu biletu na stronie internetowej, by móc uczestniczyć w wydarzeniu. [20] "

Replace-Text `
    "Jako organizator wydarzenia chce mieć wgląd do listy wolontariuszy, aby móc przypisywać im zadania do wykonania." `
    "Jako organizator wydarzenia chce mieć wgląd do listy wolontariuszy, aby móc przypisywać im zadania do wykonania. [5]"

Replace-Text `
    "Jako uzytkownik chce mieć mozliwosc, przeczytania opisu wydarzenia, by moc poznac wiecej szczegolow" `
    "Jako uzytkownik chce mieć możliwość, przeczytania opisu wydarzenia, by moc poznac wiecej szczegolow. [3]"

Replace-Text `
    "Jako kierownik chce mieć mozliwosc zmienienia statusu dla pracownika, w kalendarzu" `
    "Jako kierownik chce mieć możliwość zmienienia statusu dla pracownika, w kalendarzu [3]"

Replace-Text `
    "Jako kierownik chce mieć mozliwosc ustalanie zadan dla wolontariuszy" `
    "Jako kierownik chce mieć możliwość ustalania zadań dla wolontariuszy. [13]"

Replace-Text `
    "Jako kierownik chce mieć możliwość zapisu kontaktów do innych firm" `
    "Jako kierownik chce mieć możliwość zapisu kontaktów do innych firm.[2]"

# Insert a new, empty paragraph right after the "zapisu kontaktów" paragraph
# (and before the pre-existing blank paragraph at the end of the document).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*zapisu kontaktów do innych firm*") {
        $para.Range.InsertParagraphAfter()
        break
    }
}
